$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Status column (E) values
$ws.Range("E2").Value = "Em Andamento"
$ws.Range("E5").Value = "Em Andamento"
$ws.Range("E6").Value = "Em Andamento"
$ws.Range("E7").Value = "Em Andamento"
$ws.Range("E8").Value = "Em Andamento"

# Update Porcentagem column (F) values
$ws.Range("F2").Value = 0.75
$ws.Range("F5").Value = 0.15
$ws.Range("F6").Value = 0.35
$ws.Range("F7").Value = 0.25
$ws.Range("F8").Value = 1
$ws.Range("F12").Value = 0.25

# Update the date in I4
$ws.Range("I4").Value = "23/09/2022"

# Update the active selection / scroll position to match the saved view
$ws.Range("C5").Select() | Out-Null
